{"js": "const replacements = [\n  [\"2023-07-29 Saturday\", \"2023-07-30 Sunday\"],\n  [\"92\u00f73=\", \"43\u00f74=\"],\n  [\"92\u00f72=\", \"32\u00f76=\"],\n  [\"64\u00f75=\", \"93\u00f77=\"],\n  [\"76\u00f74=\", \"26\u00f73=\"],\n  [\"45\u00f72=\", \"62\u00f79=\"],\n  [\"31\u00f77=\", \"71\u00f72=\"],\n  [\"35\u00f78=\", \"91\u00f79=\"],\n  [\"93\u00f75=\", \"45\u00f73=\"],\n  [\"83\u00f77=\", \"48\u00f79=\"],\n  [\"25\u00f75=\", \"45\u00f77=\"],\n  [\"58\u00f74=\", \"48\u00f75=\"],\n  [\"65\u00f72=\", \"74\u00f76=\"],\n  [\"92\u00f78=\", \"77\u00f79=\"],\n  [\"47\u00f75=\", \"47\u00f76=\"],\n  [\"83\u00f79=\", \"63\u00f75=\"],\n  [\"65\u00f73=\", \"55\u00f77=\"],\n  [\"73\u00f76=\", \"35\u00f79=\"],\n  [\"19\u00f79=\", \"64\u00f73=\"],\n  [\"25\u00f73=\", \"82\u00f74=\"],\n  [\"21\u00f79=\", \"30\u00f75=\"],\n  [\"47\u00f79=\", \"79\u00f78=\"],\n  [\"87\u00f79=\", \"94\u00f72=\"],\n  [\"72\u00f75=\", \"23\u00f79=\"],\n  [\"60\u00f73=\", \"46\u00f76=\"],\n  [\"15\u00f74=\", \"42\u00f78=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-07-29 Saturday\", \"2023-07-30 Sunday\"),\n    @(\"92\u00f73=\", \"43\u00f74=\"),\n    @(\"92\u00f72=\", \"32\u00f76=\"),\n    @(\"64\u00f75=\", \"93\u00f77=\"),\n    @(\"76\u00f74=\", \"26\u00f73=\"),\n    @(\"45\u00f72=\", \"62\u00f79=\"),\n    @(\"31\u00f77=\", \"71\u00f72=\"),\n    @(\"35\u00f78=\", \"91\u00f79=\"),\n    @(\"93\u00f75=\", \"45\u00f73=\"),\n    @(\"83\u00f77=\", \"48\u00f79=\"),\n    @(\"25\u00f75=\", \"45\u00f77=\"),\n    @(\"58\u00f74=\", \"48\u00f75=\"),\n    @(\"65\u00f72=\", \"74\u00f76=\"),\n    @(\"92\u00f78=\", \"77\u00f79=\"),\n    @(\"47\u00f75=\", \"47\u00f76=\"),\n    @(\"83\u00f79=\", \"63\u00f75=\"),\n    @(\"65\u00f73=\", \"55\u00f77=\"),\n    @(\"73\u00f76=\", \"35\u00f79=\"),\n    @(\"19\u00f79=\", \"64\u00f73=\"),\n    @(\"25\u00f73=\", \"82\u00f74=\"),\n    @(\"21\u00f79=\", \"30\u00f75=\"),\n    @(\"47\u00f79=\", \"79\u00f78=\"),\n    @(\"87\u00f79=\", \"94\u00f72=\"),\n    @(\"72\u00f75=\", \"23\u00f79=\"),\n    @(\"60\u00f73=\", \"46\u00f76=\"),\n    @(\"15\u00f74=\", \"42\u00f78=\"),\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($pair[0], $false, $false, $false, $false, $false, $true, $true, $null, $pair[1], 2) | Out-Null\n}\n"}
